# 4bit unity frame encoding
# Fill in rows 4-7 with the new trial configurations, update the H column
# (notes) width, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: short, head-fixed-style trial (no special style) ---
$ws.Range("A4").Value = "[-90, 90]"
$ws.Range("B4").Value = "[2]"
$ws.Range("C4").Value = "[0.04444]"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = "short"
$ws.Range("C4:D4").Style = "Normal"

# --- Row 5 ---
$ws.Range("A5").Value = "[-90, 90]"
$ws.Range("B5").Value = "[2]"
$ws.Range("C5").Value = "[0.04444]"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 1
$ws.Range("H5").Value = "short2"
$ws.Range("C5:D5").Style = "Normal"

# --- Row 6 (keeps the existing C/D styling, like row 3) ---
$ws.Range("A6").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"
$ws.Range("B6").Value = "[2]"
$ws.Range("C6").Value = "[0.04]"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 6
$ws.Range("H6").Value = "head fixed"

# --- Row 7 (keeps the existing C/D styling, like row 3) ---
$ws.Range("A7").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"
$ws.Range("B7").Value = "[2]"
$ws.Range("C7").Value = "[0.04]"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("H7").Value = "freely moving"

# --- Widen the notes column to fit the new longer values ---
$ws.Columns(8).ColumnWidth = 13

# --- Move the active selection ---
$ws.Range("Q11").Select()
